$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block right under the header (rows 2-7: THAIS, TIAGO, JOAO, SABRINA,
# THIAGO, RODRIGO) is replaced by a refreshed export with one extra row
# (rows 2-8): some of the old accounts keep their prior balance, THAIS gets a
# new balance, three brand-new accounts show up (THEOMAR, CAIO, ANDRE), and
# two old accounts (TIAGO, SABRINA) drop off the list.
# Insert one fresh row so the 6-row block becomes 7 rows, then stamp the final
# account/name/balance values into rows 2-8.
$ws.Range("A8").EntireRow.Insert() | Out-Null

$newRows = @(
    @("004381328", "JOAO", 35538.54),
    @("005064129", "THIAGO", 20354.42),
    @("004231509", "THEOMAR", 953.1),
    @("004392159", "RODRIGO", 900.21),
    @("004512434", "CAIO", 720),
    @("005395948", "THAIS", 424.86),
    @("005003629", "ANDRE", 381.24)
)

$r = 2
foreach ($row in $newRows) {
    # Account numbers are zero-padded strings ("004381328"), so force the
    # cell to text format before assigning - otherwise Excel would coerce
    # the numeric-looking text to a number and drop the leading zeros.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
